$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E5").Value = 6.66
$ws.Range("K5").Value = 4.4029
$ws.Range("AA5").Value = 0.5893695799999999
$ws.Range("AB5").Value = 476
$ws.Range("AY5").Value = 11
$ws.Range("AZ5").Value = 35.483870967742
$ws.Range("BA5").Value = 48
$ws.Range("BB5").Value = 48.484848484848
$ws.Range("BC5").Value = 37
$ws.Range("BD5").Value = 48.684210526316
$ws.Range("BQ5").Value = 133.2
$ws.Range("CS5").Value = 31
$ws.Range("E6").Value = 6.8818181818182
$ws.Range("K6").Value = 2.2361
$ws.Range("AA6").Value = 2.53990959
$ws.Range("AB6").Value = 781
$ws.Range("AE6").Value = 343
$ws.Range("AF6").Value = 77.42663656884901
$ws.Range("AI6").Value = 235
$ws.Range("AJ6").Value = 144
$ws.Range("AY6").Value = 18
$ws.Range("AZ6").Value = 60
$ws.Range("BA6").Value = 81
$ws.Range("BB6").Value = 49.390243902439
$ws.Range("BC6").Value = 70
$ws.Range("BD6").Value = 52.631578947368
$ws.Range("BG6").Value = 226
$ws.Range("BQ6").Value = 151.4
$ws.Range("BT6").Value = 100
$ws.Range("CS6").Value = 30
$ws.Range("AA7").Value = 1.613539
$ws.Range("AB7").Value = 588
$ws.Range("AL7").Value = 50
$ws.Range("BG7").Value = 124
$ws.Range("CZ7").Value = 18
$ws.Range("E8").Value = 6.5705882352941
$ws.Range("K8").Value = 1.033
$ws.Range("AA8").Value = 1.08096045
$ws.Range("AB8").Value = 629
$ws.Range("AE8").Value = 321
$ws.Range("AF8").Value = 82.519280205656
$ws.Range("AG8").Value = 389
$ws.Range("AI8").Value = 167
$ws.Range("AJ8").Value = 85
$ws.Range("AN8").Value = 15.384615384615
$ws.Range("AQ8").Value = 5
$ws.Range("AR8").Value = 43
$ws.Range("BQ8").Value = 111.7
$ws.Range("BT8").Value = 68
$ws.Range("CC8").Value = 21
$ws.Range("CD8").Value = 11
$ws.Range("CT8").Value = 26
$ws.Range("DG8").Value = 207
$ws.Range("E9").Value = 6.7636363636364
$ws.Range("K9").Value = 4.4906
$ws.Range("AA9").Value = 6.4410427
$ws.Range("AR9").Value = 97
$ws.Range("BQ9").Value = 148.8
$ws.Range("DA9").Value = 17
$ws.Range("DB9").Value = 56.666666666667
$ws.Range("K10").Value = 3.5465
$ws.Range("AA10").Value = 2.2440736
$ws.Range("AB10").Value = 1042
$ws.Range("AE10").Value = 595
$ws.Range("AF10").Value = 76.478149100257
$ws.Range("AG10").Value = 778
$ws.Range("AH10").Value = 317
$ws.Range("DF10").Value = 362
$ws.Range("E11").Value = 6.725
$ws.Range("K11").Value = 0.0592
$ws.Range("AA11").Value = 0.63423926
$ws.Range("AB11").Value = 399
$ws.Range("AF11").Value = 79.741379310345
$ws.Range("AG11").Value = 232
$ws.Range("BA11").Value = 41
$ws.Range("BB11").Value = 65.079365079365
$ws.Range("BE11").Value = 10
$ws.Range("BF11").Value = 76.92307692307701
$ws.Range("BG11").Value = 76
$ws.Range("BQ11").Value = 53.8
$ws.Range("BT11").Value = 47
$ws.Range("CU11").Value = 22
$ws.Range("CV11").Value = 3
$ws.Range("DG11").Value = 128
$ws.Range("AA12").Value = 0.46281827
$ws.Range("AA13").Value = 1.30527372
$ws.Range("AB13").Value = 1341
$ws.Range("AE13").Value = 902
$ws.Range("AF13").Value = 86.897880539499
$ws.Range("AG13").Value = 1038
$ws.Range("AH13").Value = 493
$ws.Range("AR13").Value = 95
$ws.Range("CC13").Value = 63
$ws.Range("CD13").Value = 42
$ws.Range("DF13").Value = 549
$ws.Range("E15").Value = 6.552380952381
$ws.Range("AB15").Value = 1280
$ws.Range("BG15").Value = 291
$ws.Range("BQ15").Value = 137.6
$ws.Range("E16").Value = 6.447619047619
$ws.Range("AA16").Value = 0.43180164
$ws.Range("BQ16").Value = 135.4
$ws.Range("E18").Value = 7.01
$ws.Range("AA18").Value = 0.8904564
$ws.Range("BQ18").Value = 140.2
$ws.Range("E19").Value = 6.6923076923077
$ws.Range("AA19").Value = 0.68696183
$ws.Range("AE19").Value = 188
$ws.Range("AF19").Value = 78.661087866109
$ws.Range("AI19").Value = 107
$ws.Range("AJ19").Value = 63
$ws.Range("AL19").Value = 50
$ws.Range("BG19").Value = 81
$ws.Range("BQ19").Value = 87
$ws.Range("BT19").Value = 51
$ws.Range("CD19").Value = 6
$ws.Range("CZ19").Value = 10
$ws.Range("AA22").Value = 0.03518422
$ws.Range("AB22").Value = 260
$ws.Range("AE22").Value = 167
$ws.Range("AF22").Value = 85.641025641026
$ws.Range("AG22").Value = 195
$ws.Range("AH22").Value = 127
$ws.Range("DF22").Value = 139
$ws.Range("AA24").Value = 0.02623273
$ws.Range("AK24").Value = 92
$ws.Range("AL24").Value = 48.167539267016
$ws.Range("CZ24").Value = 191
$ws.Range("DL24").Value = -0.2269
